# NAP2 Financial data verif ganti time out
#
# Updates the raw input figures on the "Gross Yield (CF)" sheet of the
# Regular Fixed CF simulation workbook. Only the source cells themselves
# are touched (B2, B5, B6, E6, H6, E7, H7, B8, D24) -- none of the
# dependent formulas elsewhere on the sheet were recalculated/re-saved by
# the author, so we switch the workbook to manual calculation before
# writing the new values to avoid rippling a recalculation into the
# formula cells that reference them (B7, the IRR block, etc.).
$xlCalculationManual = -4135
$excel.Calculation = $xlCalculationManual

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Gross Yield (CF)")

$ws.Range("B2").Value = 350620000      # was 1.22605E8
$ws.Range("B5").Value = 1002000000     # was 6.11E8
$ws.Range("B6").Value = 350620000      # was 1.2262E8
$ws.Range("E6").Value = 39270500       # was 3.3151E7
$ws.Range("H6").Value = 29325250       # was 2.476125E7
$ws.Range("E7").Value = 0              # was 4.3455345E7
$ws.Range("H7").Value = 0              # was 4.44542E7
$ws.Range("B8").Value = 0.14685478     # was 0.1468207
$ws.Range("D24").Value = 36437000      # was 3.08145E7
